# Insert new column D (shifts existing D:K data to E:L) and populate the
# newly inserted most-recent-period column, along with the handful of
# restated prior-period cells (old D/E values corrected on re-publish).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUL")

# Insert a new column before column D; this shifts D:K -> E:L automatically,
# matching formulas/values, and updates the sheet dimension.
$ws.Columns.Item(4).Insert()

# Copy number formatting from the (now-shifted) column E into the new column D
# so the new "most recent period" column matches the date/number formats used
# by the other period columns.
$ws.Columns.Item(5).Copy()
$ws.Columns.Item(4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# (row, column, value) triples: new column-D figures for the newly added
# reporting period, plus small restatements to the two prior periods that
# shifted into columns E/F.
$data = @(
    @(7, 4, 43435),
    @(8, 4, 3041000),
    @(9, 4, 2198400),
    @(9, 5, 1691200),
    @(9, 6, 1488800),
    @(10, 4, 842600),
    @(10, 5, 614800),
    @(10, 6, 605800),
    @(12, 4, 35500),
    @(13, 4, 0),
    @(14, 4, 8900),
    @(14, 6, -300),
    @(15, 4, 0),
    @(17, 4, 2786200),
    @(17, 5, 2178000),
    @(17, 6, 1896300),
    @(18, 4, 254800),
    @(18, 5, 128000),
    @(18, 6, 198400),
    @(20, 4, 13000),
    @(21, 4, 412800),
    @(21, 5, 191600),
    @(21, 6, 268500),
    @(22, 4, 111000),
    @(23, 4, 156700),
    @(23, 5, 60600),
    @(23, 6, 163400),
    @(24, 4, 31100),
    @(24, 5, 9800),
    @(24, 6, 48900),
    @(25, 4, 0),
    @(26, 4, 125600),
    @(26, 5, 50800),
    @(26, 6, 114500),
    @(27, 4, 133700),
    @(27, 5, 59400),
    @(27, 6, 121700),
    @(28, 4, 0),
    @(29, 4, 37500),
    @(30, 4, 0),
    @(31, 4, 0),
    @(32, 4, -13000),
    @(33, 4, 171200),
    @(33, 5, 59400),
    @(33, 6, 121700),
    @(34, 4, 0),
    @(35, 4, 171200),
    @(35, 5, 59400),
    @(35, 6, 121700),
    @(38, 4, 43435),
    @(41, 4, 150800),
    @(42, 4, 0),
    @(43, 4, 512300),
    @(44, 4, 355600),
    @(44, 5, 359000),
    @(45, 4, 69100),
    @(46, 4, 1087700),
    @(46, 5, 1157600),
    @(47, 4, 108500),
    @(47, 5, 77800),
    @(48, 4, 636500),
    @(49, 4, 2213300),
    @(50, 4, 0),
    @(51, 4, 0),
    @(52, 4, 129200),
    @(52, 5, 129200),
    @(53, 4, 0),
    @(54, 4, 4175300),
    @(54, 5, 4373200),
    @(57, 4, 273400),
    @(58, 4, 106000),
    @(59, 4, 166800),
    @(60, 4, 546100),
    @(61, 4, 2141500),
    @(62, 4, 335400),
    @(62, 5, 417600),
    @(63, 4, 0),
    @(64, 4, 0),
    @(65, 4, 0),
    @(66, 4, 3023500),
    @(66, 5, 3321800),
    @(68, 4, 0),
    @(69, 4, 0),
    @(70, 4, 0),
    @(71, 4, 0),
    @(72, 4, 1285200),
    @(72, 5, 1127000),
    @(73, 4, 0),
    @(74, 4, 0),
    @(75, 4, 0),
    @(76, 4, 1151800),
    @(76, 5, 1051400),
    @(77, 4, 0),
    @(80, 4, 43435),
    @(81, 4, 171200),
    @(81, 5, 59400),
    @(81, 6, 121700),
    @(83, 4, 145100),
    @(84, 4, 0),
    @(85, 4, 0),
    @(86, 4, 0),
    @(87, 4, 0),
    @(88, 4, 0),
    @(89, 4, 253300),
    @(91, 4, -68300),
    @(92, 4, 0),
    @(93, 4, 0),
    @(94, 4, -61800),
    @(96, 4, -31100),
    @(97, 4, 0),
    @(98, 4, 0),
    @(99, 4, 0),
    @(100, 4, -228600),
    @(101, 4, -6500),
    @(102, 4, -43600)
)

foreach ($item in $data) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws.Cells.Item($r, $c).Value = $v
}
